$d = $word.ActiveDocument

function Strip-TrailingControlChars($s) {
    # Paragraph.Range.Text always carries the trailing paragraph mark (\r,
    # code 13) and, for the last paragraph inside a table cell, an
    # additional cell-end marker (\a, code 7). Strip any such trailing
    # control characters (< 0x20) so we can compare against plain text.
    $end = $s.Length
    while ($end -gt 0 -and [int][char]$s.Substring($end - 1, 1) -lt 32) {
        $end--
    }
    return $s.Substring(0, $end)
}

function Find-ParagraphByText($doc, $exactText) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $txt = Strip-TrailingControlChars $p.Range.Text
        if ($txt -eq $exactText) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Edit 1: merge field "<<cs_{defendant1.individual}>>" was referencing the
# wrong / incomplete property name. Fix it to use the real boolean property
# "isIndividual" -> "<<cs_{defendant1.isIndividual }>>"
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphByText $d "<<cs_{defendant1.individual}>>"
if ($p1 -eq $null) {
    throw "Could not find target paragraph for edit 1 (cs_{defendant1.individual})"
}

$p1Start = $p1.Range.Start
$p1Text = $p1.Range.Text
$relIdx = $p1Text.IndexOf("individual}")
if ($relIdx -lt 0) {
    throw "Could not locate 'individual}' inside target paragraph for edit 1"
}
$iAbs = $p1Start + $relIdx

# Replace the tail "ndividual}>>" (everything after the leading "i") with "}>>"
# first, so the offset of the leading "i" run stays valid.
$tailLen = "ndividual}>>".Length
$rTail = $d.Range($iAbs + 1, $iAbs + 1 + $tailLen)
if ($rTail.Text -ne "ndividual}>>") {
    throw "Unexpected text in tail range for edit 1: $($rTail.Text)"
}
$rTail.Text = "}>>"

# Replace the leading "i" with "isIndividual " (adds the missing suffix plus
# a separating space before the closing "}>>").
$rHead = $d.Range($iAbs, $iAbs + 1)
if ($rHead.Text -ne "i") {
    throw "Unexpected text in head range for edit 1: $($rHead.Text)"
}
$rHead.Text = "isIndividual "

# ---------------------------------------------------------------------------
# Edit 2: merge field "<<cr_{defendant1.correspondenceAddress != null}>>" had
# a stray space before "!=" that doesn't match the other, correctly-formatted
# conditions elsewhere in the template (e.g. "correspondenceAddress!= null").
# Remove that single space.
# ---------------------------------------------------------------------------
$p2 = Find-ParagraphByText $d "<<cr_{defendant1.correspondenceAddress != null}>>"
if ($p2 -eq $null) {
    throw "Could not find target paragraph for edit 2 (cr_{defendant1.correspondenceAddress != null})"
}

$p2Start = $p2.Range.Start
$p2Text = $p2.Range.Text
$relIdx2 = $p2Text.IndexOf("correspondenceAddress ")
if ($relIdx2 -lt 0) {
    throw "Could not locate 'correspondenceAddress ' inside target paragraph for edit 2"
}
$spaceAbs = $p2Start + $relIdx2 + "correspondenceAddress".Length

$rSpace = $d.Range($spaceAbs, $spaceAbs + 1)
if ($rSpace.Text -ne " ") {
    throw "Unexpected text in space range for edit 2: $($rSpace.Text)"
}
$rSpace.Text = ""

Write-Output "Edit1 -> $(Strip-TrailingControlChars $p1.Range.Text)"
Write-Output "Edit2 -> $(Strip-TrailingControlChars $p2.Range.Text)"
